# accesskarteringinfo.xlsx update
# commit message: "added crs check; filled in shapefilekarteringinfo"
#
# Fills in several previously-placeholder "elm_id_name" (column D) values,
# updates two "path_csvs"/"path_shapes" (columns C/E) entries for the
# Vlieland_2013 and "NM vegetatiekartering RuitenAa2020" rows (the latter's
# source folder was renamed to "... RuitenAa2020 edited"), and widens
# columns C and D to fit the now-longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Drenthe extra / NM_DCR_2013 ---------------------------------
# elm_id_name: "OBJECTID (denk ik)" -> "OBJECTID"
$ws.Range("D3").Value = "OBJECTID"

# --- Row 4: Drenthe / NM_Westerbork2017 ---------------------------------
# elm_id_name: filled in with the actual field name "ELMID"
$ws.Range("D4").Value = "ELMID"

# --- Row 14: Friesland / Vlieland_2013 ----------------------------------
# path_csvs updated to the correct digital-standard database name
$ws.Range("C14").Value = "./FR/Duinen_Vlieland/Vegetatiekartering 2013/890_Vlieland2013"

# --- Row 23: Groningen / NM vegetatiekartering RuitenAa2020 -------------
# source data folder renamed to "... RuitenAa2020 edited"
$ws.Range("C23").Value = "./GR/NM vegetatiekartering RuitenAa2020 edited/digi standaard_RuitenA_2020/digi standaard_RuitenA_2020"
$ws.Range("E23").Value = "./GR/NM vegetatiekartering RuitenAa2020 edited/vegkart_RuitenA_2020/vegkart_RuitenA_2020.shp"

# --- Column widths: widen C and D to fit the longer path/id text --------
$ws.Columns.Item(3).ColumnWidth = 126.83
$ws.Columns.Item(4).ColumnWidth = 63.17

# --- Selection / scroll position, matching the saved view ---------------
$ws.Range("E24").Select()
$excel.ActiveWindow.ScrollColumn = 3
